# Certificate model: remove the RG (ID card) reference from the
# "portador do RG nº {{RG}} e CPF nº {{CPF}}" sentence, leaving just
# "portador do CPF nº {{CPF}}" — per commit "remove RG to user and
# models certificate".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)          # "Rectangle 5" - the certificate body text box
$tr = $sh.TextFrame.TextRange

# The sentence lives in paragraph 1 of this text box, originally built out of
# these consecutive runs (1-based character offsets into the full TextRange):
#   (28,18) "portador do RG nº "
#   (46,2)  "{{"
#   (48,2)  "RG"
#   (50,3)  "}} "
#   (53,9)  "e CPF nº "
#   (62,9)  "{{CPF}}, "            <- keep this one untouched
#
# Target wording: "portador " + "do CPF " + "nº " + "{{CPF}}, "
#
# Edits are applied from the right-most offset to the left-most so that
# earlier (lower) offsets stay valid while later ones are being rewritten.

$tr.Characters(53, 9).Text = "n" + [char]0x00BA + " "   # "e CPF nº "  -> "nº "
$tr.Characters(50, 3).Text = ""                          # "}} "       -> ""
$tr.Characters(48, 2).Text = "do CPF "                    # "RG"        -> "do CPF "
$tr.Characters(46, 2).Text = ""                           # "{{"        -> ""
$tr.Characters(37, 9).Text = ""                           # "do RG nº " -> ""
$tr.Characters(28, 9).Text = "portador "                  # "portador " -> "portador " (splits the run)
